$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update issue statuses (before the list gets re-sorted by the sheet's sort order)
# Issue_010 "Create new assert UI takes too long to load": Open -> In Work
$ws.Range("D5").Value = "In Work"
# Issue_031 "Use vim's for templates to reduce overhead": In work -> Solved
$ws.Range("D2").Value = "Solved"
# Issue_032 "Bin the inline insert": Open -> Solved
$ws.Range("D4").Value = "Solved"
# Issue_034 "Add a default description and other parameters...": In work -> Solved
$ws.Range("D3").Value = "Solved"
# Issue_033 "Rearchitect Build Evaluation String to use the test class": Open -> Solved
$ws.Range("D6").Value = "Solved"

# Re-apply the worksheet's existing sort order (Status asc, Priority desc, ID asc)
$dataRange = $ws.Range("A2:F35")
$keyStatus = $ws.Range("D2:D35")
$keyPriority = $ws.Range("C2:C35")
$keyId = $ws.Range("A2:A35")
$dataRange.Sort($keyStatus, 1, $keyPriority, [Type]::Missing, 2, $keyId, 1, 2)

# Move the active selection to D3, matching the saved view state
$ws.Range("D3").Select()
